$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate pair_kind ("generic") for the four practice pair rows.
$ws.Range("J2:J5").Value = "generic"

# New "stim details" block appended below the existing table.
$ws.Range("A27").Value = "stim details"

$ws.Range("A28").Value = "month"
$ws.Range("B28").Value = "word_type"
$ws.Range("C28").Value = "need_audio"
$ws.Range("D28").Value = "need_image"
$ws.Range("E28").Value = "word"
$ws.Range("F28").Value = "count"
$ws.Range("G28").Value = "find images"

$videoRows = 29, 30, 31, 32
$videoMonths = 6, 6, 7, 7
for ($i = 0; $i -lt $videoRows.Length; $i++) {
    $r = $videoRows[$i]
    $ws.Cells.Item($r, 1).Value = $videoMonths[$i]
    $ws.Cells.Item($r, 2).Value = "video"
}

$audioRows = 33, 34, 35, 36
$audioMonths = 6, 6, 7, 7
for ($i = 0; $i -lt $audioRows.Length; $i++) {
    $r = $audioRows[$i]
    $ws.Cells.Item($r, 1).Value = $audioMonths[$i]
    $ws.Cells.Item($r, 2).Value = "audio"
}
